# Adding black bear (Ursus americanus) as a new watchlist-species row.
# It is inserted immediately above "Uvularia sessilifolia" (row 89), in the
# "rare native" / "P" (present in ANP) block, shifting every following row
# down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 89, pushing existing rows 89+ down to 90+.
$ws.Rows("89:89").Insert()

# Populate the new row with the black bear entry.
$ws.Cells.Item(89, 1).Value = "Ursus americanus"
$ws.Cells.Item(89, 2).Value = "rare native"
$ws.Cells.Item(89, 3).Value = "P"

# Restore the view state (best effort - scroll position / zoom / selection).
$win = $excel.ActiveWindow
$win.ScrollRow = 69
$win.ScrollColumn = 1
$win.Zoom = 117
$ws.Range("E88").Select()
